$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.048.27"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "3.005.09"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.30"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.78"
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.03"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0855"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.02"
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").Value = "3.479.21"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.59"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").Value = "2.989.05"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "52.042.24"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.45"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.52"
$ws.Range("E21").Value = "  -3.83%  "
$ws.Range("D22").Value = "0.0₃0969"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.07"
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.75"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.72"
$ws.Range("E25").Value = "  -3.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.179"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.94"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.110"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.34"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.16"
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.09"
$ws.Range("E33").Value = "  -7.63%  "
$ws.Range("E34").Value = "  +13.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.29"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0433"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.48"
$ws.Range("E41").Value = "  -5.51%  "
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.23"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.70"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "2.125.56"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("E48").Value = "  -7.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.244"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.904"
$ws.Range("E51").Value = "  -0.55%  "
